$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new columns E and F ---
# Copy formatting (and content) from D1 (bold/centered header style) onto
# E1 and F1, then overwrite their text values with the new header labels.
$ws.Range("D1").Copy($ws.Range("E1"))
$ws.Range("D1").Copy($ws.Range("F1"))
$ws.Range("E1").Value = "shock_raw"
$ws.Range("F1").Value = "extreme_level_raw"

# --- shock_raw (column E): the raw numeric form of the "shock" column ---
$ws.Range("E2").Value = 0.2
$ws.Range("E3").Value = 0.04347826086956519
$ws.Range("E4").Value = 0.02439024390243905
$ws.Range("E5").Value = 0.01694915254237284
$ws.Range("E6").Value = 0.01298701298701288
$ws.Range("E7").Value = 0.01052631578947372
$ws.Range("E8").Value = 0.008849557522123908
$ws.Range("E9").Value = 0.007633587786259444
$ws.Range("E10").Value = 0.006711409395973256
$ws.Range("E11").Value = 0.07784431137724557
$ws.Range("E12").Value = 13
$ws.Range("E13").Value = 13
$ws.Range("E14").Value = 13
$ws.Range("E15").Value = 13
$ws.Range("E16").Value = 0.01118870236745795
$ws.Range("E17").Value = 13
$ws.Range("E18").Value = 13
$ws.Range("E19").Value = 13
$ws.Range("E20").Value = 0.002386634844868674
$ws.Range("E21").Value = 13

# --- extreme_level_raw (column F): the raw numeric form of "extreme_level" ---
# Rows with a ppts/bps annotation or numeric ratio in D get the matching
# raw number; rows where D is blank get a blank (empty-text) F cell too,
# matching D's own blank-but-present cell.
$blankRows = @(2,3,4,5,6,7,8,9,10,11,20)
foreach ($r in $blankRows) {
    $cell = $ws.Range("F$r")
    $cell.Value = "'"
    $cell.Style = "Normal"
}

$ws.Range("F12").Value = 198
$ws.Range("F13").Value = 216
$ws.Range("F14").Value = 234
$ws.Range("F15").Value = 252
$ws.Range("F16").Value = 0.01157730348796671
$ws.Range("F17").Value = 378
$ws.Range("F18").Value = 396
$ws.Range("F19").Value = 414
$ws.Range("F21").Value = 450

Write-Output "shock_raw / extreme_level_raw columns written"
